$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from serial date 45175 to 45183 (2023-09-06 -> 2023-09-14)
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45183
}
